$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -152.5562
$ws.Range("B2").Value = -152.3113

$ws.Range("A3").Value = 58.885
$ws.Range("B3").Value = 59.0341

$ws.Range("A4").Value = -150.2082
$ws.Range("B4").Value = -150.4509

$ws.Range("A5").Value = 60.3313
$ws.Range("B5").Value = 60.18
